$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the existing rows: the label set grew from 3 classes
# (COMMA/PERIOD/O) to 5 (O/,/./?/<PAD>), and the per-class attention
# scores were refreshed. ---
$ws.Range("B2").Value = "O"
$ws.Range("C2").Value = 0.8727850976828715
$ws.Range("D2").Value = 0.7860065466448445
$ws.Range("E2").Value = 0.9811031664964249

$ws.Range("B3").Value = ","
$ws.Range("C3").Value = 0.5094339622641509
$ws.Range("D3").Value = 0.6192660550458715
$ws.Range("E3").Value = 0.4326923076923077

$ws.Range("B4").Value = "."
$ws.Range("C4").Value = 0.2780082987551867
$ws.Range("D4").Value = 0.5317460317460317
$ws.Range("E4").Value = 0.1882022471910112

$ws.Range("B5").Value = "?"
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0

$ws.Range("B6").Value = "<PAD>"
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0

# --- Append the Macro / Micro average rows (previously rows 5/6, now
# pushed down to rows 7/8 by the expanded label set). Copy the number
# formatting/style used by column A on the preceding rows first, then
# fill in the actual values. ---
$ws.Range("A6").Copy()
$ws.Range("A7:A8").PasteSpecial(-4122)

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "Macro"
$ws.Range("C7").Value = 0.3507301606914083
$ws.Range("D7").Value = 0.3874037266873495
$ws.Range("E7").Value = 0.3203995442759487

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "Micro"
$ws.Range("C8").Value = 0.2324510932105869
$ws.Range("D8").Value = 0.2329873125720877
$ws.Range("E8").Value = 0.2319173363949483
